# Commit: "Fixing Bulk Operation Template"
# Rename the worksheet tab from "Create Item Category" to "Create Category"
# (xl/workbook.xml: <sheet name="Create Item Category" .../> -> <sheet name="Create Category" .../>)

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Create Item Category") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Name = "Create Category"
